$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H15").Value = 0
